# Apply "last minute updates" to the first paragraph of the document:
#   - add a paragraph border (space-only, no line) of 5 on all four sides
#   - change the left indent from 120 (6pt) to 225 (11.25pt) twentieths of a point
#   - update the bookmark-style ID text and drop the trailing whitespace run

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)

# Replace the ID placeholder text. Matching the original text together with
# its trailing space (which lived in a second, separately-formatted run)
# and replacing it with the new text (no trailing space) collapses the
# paragraph back down to a single run, just like in the target document.
$p1.Range.Find.Execute(
    "**ID__AFFARS_pgi_5322_topic_3__ID** ", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5322_101_1_90__ID**", 2
) | Out-Null

$pf = $d.Paragraphs.Item(1).Range.ParagraphFormat

# Update the left indent: 120 twips (6 pt) -> 225 twips (11.25 pt)
$pf.LeftIndent = 11.25

# Add a paragraph border on all sides with only a "space" (distance from
# text) of 5, and no line style/width/color -- matching <w:pBdr> with only
# w:space="5" on each edge.
$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5
